$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.970.71"
$ws.Range("E2").Value = "  -1.27%  "

$ws.Range("D3").Value = "2.356.06"
$ws.Range("E3").Value = "  -1.40%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.32"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.88"
$ws.Range("E6").Value = "  -1.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.540"
$ws.Range("E8").Value = "  -2.41%  "

$ws.Range("D9").Value = "2.369.03"
$ws.Range("E9").Value = "  -1.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0960"
$ws.Range("E10").Value = "  -1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -0.69%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.79"
$ws.Range("E12").Value = "  +2.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.322"
$ws.Range("E13").Value = "  -4.36%  "

$ws.Range("D14").Value = "2.779.40"
$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.78"
$ws.Range("E15").Value = "  -0.26%  "

$ws.Range("D16").Value = "55.964.04"
$ws.Range("E16").Value = "  -1.25%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000131"
$ws.Range("E17").Value = "  -1.39%  "

$ws.Range("D18").Value = "2.391.67"
$ws.Range("E18").Value = "  +0.99%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.92"
$ws.Range("E19").Value = "  -2.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.04"
$ws.Range("E20").Value = "  -0.36%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.80"
$ws.Range("E21").Value = "  -0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("E22").Value = "  +0.41%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.28"
$ws.Range("E24").Value = "  -0.56%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  +0.14%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.147"
$ws.Range("E27").Value = "  -2.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.23"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "170.02"
$ws.Range("E29").Value = "  -1.75%  "

$ws.Range("D30").Value = "0.0₃0713"
$ws.Range("E30").Value = "  -2.36%  "

$ws.Range("E31").Value = "  -1.10%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.81"
$ws.Range("E33").Value = "  -0.75%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.08"
$ws.Range("E34").Value = "  -3.54%  "

$ws.Range("E35").Value = "  -0.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.72"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("E37").Value = "  -1.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.872"
$ws.Range("E38").Value = "  +6.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.72"
$ws.Range("E39").Value = "  -3.56%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.28"
$ws.Range("E40").Value = "  -1.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.41"
$ws.Range("E41").Value = "  -1.54%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.376"
$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.97"
$ws.Range("E44").Value = "  +2.57%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "125.71"
$ws.Range("E45").Value = "  -5.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.556"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0896"
$ws.Range("E47").Value = "  -1.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "243.22"
$ws.Range("E48").Value = "  -2.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0483"
$ws.Range("E49").Value = "  -0.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "16.98"
$ws.Range("E50").Value = "  -0.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0208"
$ws.Range("E51").Value = "  -0.97%  "
